$d = $word.ActiveDocument

function Replace-ParaXML($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End
    $r = $d.Range($start, $end)
    $full = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($full)
}

# ---------------------------------------------------------------------------
# 1) Move <w:lastRenderedPageBreak/> from the "Imagem livroLivro.png" run to
#    the preceding "Imagem input_bg.jpg" run.
# ---------------------------------------------------------------------------
$xml190 = '<w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Imagem</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> input_bg.jpg</w:t></w:r></w:p></w:body>'
Replace-ParaXML 190 $xml190

$xml191 = '<w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Imagem</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> livroLivro.png</w:t></w:r></w:p></w:body>'
Replace-ParaXML 191 $xml191

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from the "id_livro" mural variable to
#    the preceding "insere" mural variable.
# ---------------------------------------------------------------------------
$xml231 = '<w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/>' +
    '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Variável: insere = insere dados no mural</w:t></w:r></w:p></w:body>'
Replace-ParaXML 231 $xml231

$xml233 = '<w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/>' +
    '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Variável: id_livro = parametro que será inserido no mural</w:t></w:r></w:p></w:body>'
Replace-ParaXML 233 $xml233

# ---------------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from "File bootstrap-responsive.min.css"
#    to the preceding "Pacote css" run.
# ---------------------------------------------------------------------------
$xml276 = '<w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Pacote</w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> css</w:t></w:r></w:p></w:body>'
Replace-ParaXML 276 $xml276

$xml277 = '<w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">File </w:t></w:r>' +
    '<w:r><w:t>bootstrap-</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/><w:r><w:t>responsive.</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>min.css</w:t></w:r></w:p></w:body>'
Replace-ParaXML 277 $xml277

# ---------------------------------------------------------------------------
# 4) Insert a new "Variável: mensagem" bullet under
#    ExcessaoTelefoneInvalido's __construct, before "File RecebeForm.php".
# ---------------------------------------------------------------------------
$target = $d.Paragraphs.Item(163)
$target.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(163)
$insPoint = $newPara.Range.Duplicate
$insPoint.Collapse(1)
$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="27"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Variável: mensagem = contém a mensagem de ocorrência de exceção</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> na inserção do telefone</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r></w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($newParaXml)

Write-Host "Edit complete"
